# Update countries & provincias Spain
# Applies the diff: reorders Paraguay/Azerbaiyan and Bonaire/Liechtenstein
# country rows (by swapping their labels + per-row stat values), refreshes
# a handful of per-country COVID statistics, and bumps the "datos
# actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 03:39"

# --- Row 4: Estados Unidos --------------------------------------------
$ws.Range("B4").Value = 7447282
$ws.Range("C4").Value = 40929
$ws.Range("D4").Value = 4699706
$ws.Range("E4").Value = 2535836
$ws.Range("G4").Value = 955
$ws.Range("H4").Value = 211740

# --- Row 9: Peru --------------------------------------------------------
$ws.Range("B9").Value = 814829
$ws.Range("C9").Value = 3061
$ws.Range("D9").Value = 683815
$ws.Range("E9").Value = 98551
$ws.Range("G9").Value = 67
$ws.Range("H9").Value = 32463

# --- Rows 68/69: Azerbaiyan & Paraguay swap order ----------------------
# Paraguay moves above Azerbaiyan; row 68 now holds Paraguay's refreshed
# figures while row 69 keeps Azerbaiyan's previous (row 68) figures.
$ws.Range("A68").Value = "Paraguay"
$ws.Range("B68").Value = 40758
$ws.Range("C68").Value = 657
$ws.Range("D68").Value = 24449
$ws.Range("E68").Value = 15452
$ws.Range("G68").Value = 16
$ws.Range("H68").Value = 857

$ws.Range("A69").Value = "Azerbaiyan"
$ws.Range("B69").Value = 40229
$ws.Range("C69").Value = 110
$ws.Range("D69").Value = 37954
$ws.Range("E69").Value = 1684
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 591

# --- Row 82: Corea del Sur ----------------------------------------------
$ws.Range("D82").Value = 21591
$ws.Range("E82").Value = 1808

# --- Row 194: Seychelles -------------------------------------------------
$ws.Range("B194").Value = 144
$ws.Range("C194").Value = 1
$ws.Range("E194").Value = 4

# --- Rows 195/196: Liechtenstein & Bonaire swap order --------------------
# Bonaire, San Eustaquio y Saba moves above Liechtenstein; row 195 now
# holds Bonaire's refreshed figures while row 196 keeps Liechtenstein's
# previous (row 195) figures.
$ws.Range("A195").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B195").Value = 121
$ws.Range("C195").Value = 15
$ws.Range("D195").Value = 32
$ws.Range("E195").Value = 88

$ws.Range("A196").Value = "Liechtenstein"
$ws.Range("B196").Value = 118
$ws.Range("D196").Value = 113
$ws.Range("E196").Value = 4
